$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.999.49"
$ws.Range("E2").Value = "  +1.57%  "

# Row 3
$ws.Range("D3").Value = "3.526.55"
$ws.Range("E3").Value = "  +0.43%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'600.30"
$ws.Range("E5").Value = "  +1.20%  "

# Row 6
$ws.Range("D6").Value = "'183.98"
$ws.Range("E6").Value = "  +6.18%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("E8").Value = "  +0.61%  "

# Row 9
$ws.Range("D9").Value = "'0.141"
$ws.Range("E9").Value = "  +5.46%  "

# Row 10
$ws.Range("E10").Value = "  -1.75%  "

# Row 11
$ws.Range("D11").Value = "'0.444"
$ws.Range("E11").Value = "  +1.75%  "

# Row 12
$ws.Range("D12").Value = "4.139.67"
$ws.Range("E12").Value = "  +0.40%  "

# Row 13
$ws.Range("D13").Value = "'32.62"
$ws.Range("E13").Value = "  +12.80%  "

# Row 15
$ws.Range("D15").Value = "67.971.43"
$ws.Range("E15").Value = "  +1.47%  "

# Row 16
$ws.Range("E16").Value = "  +1.02%  "

# Row 17
$ws.Range("D17").Value = "3.522.69"
$ws.Range("E17").Value = "  +1.29%  "

# Row 18
$ws.Range("E18").Value = "  +1.46%  "

# Row 19
$ws.Range("D19").Value = "'14.98"
$ws.Range("E19").Value = "  +5.26%  "

# Row 20
$ws.Range("D20").Value = "'398.96"
$ws.Range("E20").Value = "  +0.89%  "

# Row 21
$ws.Range("E21").Value = "  +2.06%  "

# Row 22
$ws.Range("D22").Value = "'73.73"
$ws.Range("E22").Value = "  +0.68%  "

# Row 23
$ws.Range("E23").Value = "  +1.31%  "

# Row 24
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.22%  "

# Row 25
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "'5.70"
$ws.Range("E25").Value = "  +0.17%  "

# Row 26
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "'0.0000125"
$ws.Range("E26").Value = "  +3.17%  "

# Row 27
$ws.Range("D27").Value = "'10.76"
$ws.Range("E27").Value = "  +5.61%  "

# Row 28
$ws.Range("E28").Value = "  -0.99%  "

# Row 29
$ws.Range("E29").Value = "  -0.17%  "

# Row 30
$ws.Range("D30").Value = "'6.29"
$ws.Range("E30").Value = "  +0.43%  "

# Row 31
$ws.Range("E31").Value = "  +1.32%  "

# Row 32
$ws.Range("D32").Value = "'2.08"
$ws.Range("E32").Value = "  +1.15%  "

# Row 33
$ws.Range("D33").Value = "'24.14"
$ws.Range("E33").Value = "  +0.85%  "

# Row 34
$ws.Range("D34").Value = "'7.48"
$ws.Range("E34").Value = "  +1.32%  "

# Row 35
$ws.Range("E35").Value = "  +0.09%  "

# Row 36
$ws.Range("E36").Value = "  +2.98%  "

# Row 37
$ws.Range("D37").Value = "'163.87"
$ws.Range("E37").Value = "  +0.50%  "

# Row 38
$ws.Range("D38").Value = "'1.96"
$ws.Range("E38").Value = "  +3.00%  "

# Row 39
$ws.Range("E39").Value = "  -1.42%  "

# Row 40
$ws.Range("D40").Value = "'7.15"
$ws.Range("E40").Value = "  +3.90%  "

# Row 41
$ws.Range("E41").Value = "  +7.68%  "

# Row 42
$ws.Range("D42").Value = "'4.78"
$ws.Range("E42").Value = "  +2.26%  "

# Row 43
$ws.Range("D43").Value = "'27.12"
$ws.Range("E43").Value = "  +2.62%  "

# Row 44
$ws.Range("D44").Value = "'27.69"
$ws.Range("E44").Value = "  -0.61%  "

# Row 45
$ws.Range("D45").Value = "2.885.94"
$ws.Range("E45").Value = "  +2.92%  "

# Row 46
$ws.Range("E46").Value = "  -0.14%  "

# Row 47
$ws.Range("D47").Value = "'42.53"
$ws.Range("E47").Value = "  -0.74%  "

# Row 48
$ws.Range("D48").Value = "'352.89"
$ws.Range("E48").Value = "  +4.24%  "

# Row 49
$ws.Range("E49").Value = "  +0.56%  "

# Row 50
$ws.Range("E50").Value = "  -0.52%  "

# Row 51
$ws.Range("D51").Value = "'33.75"
$ws.Range("E51").Value = "  +1.03%  "
